# Apply crypto price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.641.25'
$ws.Range("E2").Value = '  +2.70%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.537.64'
$ws.Range("E3").Value = '  +1.45%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.73'
$ws.Range("E5").Value = '  +5.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.72'
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.618'
$ws.Range("E7").Value = '  +1.67%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.535.17'
$ws.Range("E8").Value = '  +1.70%  '
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("E10").Value = '  +5.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.73'
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.587'
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.91'
$ws.Range("E13").Value = '  +1.18%  '
$ws.Range("E14").Value = '  +2.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.107.99'
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '619.46'
$ws.Range("E16").Value = '  -4.06%  '
$ws.Range("E17").Value = '  -1.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.632.76'
$ws.Range("E18").Value = '  +2.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.526.95'
$ws.Range("E19").Value = '  +0.92%  '
$ws.Range("E20").Value = '  -1.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.43'
$ws.Range("E21").Value = '  +1.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.882'
$ws.Range("E22").Value = '  -0.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.35'
$ws.Range("E23").Value = '  -15.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.74'
$ws.Range("E24").Value = '  -1.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '96.68'
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.80'
$ws.Range("E26").Value = '  -0.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.61'
$ws.Range("E28").Value = '  -0.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.48'
$ws.Range("E29").Value = '  +2.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.11'
$ws.Range("E30").Value = '  -2.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.46'
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("E32").Value = '  -3.62%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.00'
$ws.Range("E33").Value = '  -1.19%  '
$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.30'
$ws.Range("E34").Value = '  -1.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '572.25'
$ws.Range("E35").Value = '  -7.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.101'
$ws.Range("E36").Value = '  -1.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.61'
$ws.Range("E37").Value = '  +2.25%  '
$ws.Range("E38").Value = '  +0.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '57.54'
$ws.Range("E39").Value = '  +1.76%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0470'
$ws.Range("E40").Value = '  +8.51%  '
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.142'
$ws.Range("E42").Value = '  +4.31%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.348.21'
$ws.Range("E43").Value = '  -0.18%  '
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.326'
$ws.Range("E44").Value = '  -0.85%  '
$ws.Range("E45").Value = '  +9.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '33.08'
$ws.Range("E46").Value = '  +0.51%  '
$ws.Range("E47").Value = '  +1.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.64'
$ws.Range("E48").Value = '  +3.01%  '
$ws.Range("E49").Value = '  -0.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.01'
$ws.Range("E50").Value = '  +1.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.69'
$ws.Range("E51").Value = '  -0.94%  '
